# WP excel spreadsheet - done
#
# 1. Rename the "node-f1" sheet to "nodes".
# 2. Update the view/selection on the "nodes" sheet (was node-f1): scroll to
#    A36 and select K2:K73, then make it the active/visible sheet.
# 3. Update the view/selection on the "neighbor" sheet: scroll to A82 and
#    select F2:F118 (it's no longer the active sheet).

$wb = $excel.ActiveWorkbook

$wsNodes = $wb.Worksheets.Item("node-f1")
$wsNodes.Name = "nodes"

$wsNeighbor = $wb.Worksheets.Item("neighbor")
$wsNeighbor.Activate()
$winNeighbor = $excel.ActiveWindow
$winNeighbor.ScrollRow = 82
$winNeighbor.ScrollColumn = 1
$wsNeighbor.Range("F2:F118").Select()

$wsNodes.Activate()
$winNodes = $excel.ActiveWindow
$winNodes.ScrollRow = 36
$winNodes.ScrollColumn = 1
$wsNodes.Range("K2:K73").Select()
